$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 and row 26 have effectively swapped their species-record content
# (everything except the columns that are identical between the two rows).
# Capture the "before" values for the columns that differ, then write each
# row's cells with the other row's original values.

$cols = @("A","B","D","E","F","G","H","I","J","K","Q","R","AC")

$orig25 = @{}
$orig26 = @{}
foreach ($col in $cols) {
    $orig25[$col] = $ws.Range("$col" + "25").Value2
    $orig26[$col] = $ws.Range("$col" + "26").Value2
}

# Columns whose values must stay text (they look numeric, but are stored
# as text in the sheet) - force a text number format so Excel doesn't
# silently coerce them to numbers when we write the swapped value back.
$textCols = @("I")

foreach ($col in $cols) {
    $target25 = $ws.Range("$col" + "25")
    $target26 = $ws.Range("$col" + "26")

    if ($textCols -contains $col) {
        $target25.NumberFormat = "@"
        $target26.NumberFormat = "@"
    }

    $new25 = $orig26[$col]
    $new26 = $orig25[$col]

    if (($null -eq $new25) -or ($new25 -eq "")) {
        $target25.Value2 = ""
    } else {
        $target25.Value2 = $new25
    }

    if (($null -eq $new26) -or ($new26 -eq "")) {
        $target26.Value2 = ""
    } else {
        $target26.Value2 = $new26
    }
}
